$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data (rows 18-21), matching the columns:
# A: Unnamed: 0, B: epochs, C: op_reg, D: learning_rate, E: adv_train_flag,
# F: activation, G: acc_reg, H: acc_fgsm, I: acc_pgd, J: loss_reg, K: loss_pgd,
# L: loss_fgsm, M: logdir, N: weights_path, O: tb_dir, P: Distances, Q: Norms
$rows = @(
    @(16, 20, 0, 0.003, "Regular", "<function relu at 0x11e35f9d8>", 1, 1, 1, 0.009425175376236439, 0.05570029467344284, 1, "logs/results_86.log", "weights/model_86.ckpt", "tb/86/non_robust", "(0.06946261, 0.12680605, 0.10363617)", "(117.439575, 8.788963, 6.8641105, 3.1679692)"),
    @(17, 20, 0, 0.003, "Regular", "<function relu at 0x123bb89d8>", 1, 1, 1, 0.00775056378915906, 0.1920228749513626, 1, "logs/results_87.log", "weights/model_87.ckpt", "tb/87/non_robust", "(0.18897994, 0.41342703, 0.28442824)", "(109.61839, 8.210519, 6.61972, 2.5461638)"),
    @(18, 20, 0, 0.003, "Regular", "<function relu at 0x1170ae9d8>", 1, 1, 1, 0.00669969292357564, 0.05489230528473854, 1, "logs/results_88.log", "weights/model_88.ckpt", "tb/88/non_robust", "(0.05377488, 0.06494685, 0.081667334)", "(107.84222, 7.0260696, 6.7622585, 3.9107814)"),
    @(19, 20, 0, 0.003, "Regular", "<function relu at 0x1191959d8>", 1, 1, 1, 0.005836887285113335, 0.03965282067656517, 1, "logs/results_90.log", "weights/model_90.ckpt", "tb/90/non_robust", "(0.052884344, 0.08614982, 0.049282983)", "(112.28816, 7.4369655, 6.118991, 3.4870443)")
)

$startRow = 18
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    for ($c = 0; $c -lt $data.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $data[$c]
    }
}
